# The workbook contains 16 worksheets (one per backward-elimination step),
# each holding a statsmodels OLS summary text blob in cell B2. The report
# was re-run/re-saved, so the embedded "Date:" / "Time:" stamp lines need
# to move from Sun, 05 Jan 2020 21:22:23 to Wed, 08 Jan 2020 19:07:28 on
# every sheet, leaving all other text (coefficients, stats, etc.) intact.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 05 Jan 2020"
$newDate = "Wed, 08 Jan 2020"
$oldTime = "21:22:23"
$newTime = "19:07:28"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Text
    if ($text -and $text.Contains("Date:") -and $text.Contains($oldDate)) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        $cell.Value = $updated
    }
}
